$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 48
$ws.Range("A48").Value = 'Whey 100  Pote Cookies 900GR Integralmedica'
$ws.Range("B48").Value = 10
$ws.Range("C48").Value = 184.44
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '23838'
$ws.Range("E48").Value = 161.02
$ws.Range("F48").Value = 166
$ws.Range("G48").Value = 'https://www.farmaponte.com.br/whey-100-pote-cookies-900gr-integralmedica/p'
$ws.Range("H48").NumberFormat = "@"
$ws.Range("H48").Value = '7896311766376'
$ws.Range("I48").Value = 'Integralmédica'

# Row 49
$ws.Range("A49").Value = 'Isotonico Go Drink Guarana C Acai 900gr Altetica'
$ws.Range("B49").Value = 10
$ws.Range("C49").Value = 47.67
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '25481'
$ws.Range("E49").Value = 41.61
$ws.Range("F49").Value = 42.9
$ws.Range("G49").Value = 'https://www.farmaponte.com.br/isotonico-go-drink-guarana-c-acai-900gr-altetica/p'
$ws.Range("H49").NumberFormat = "@"
$ws.Range("H49").Value = '7899621100618'
$ws.Range("I49").Value = 'Sem marca'

# Row 50
$ws.Range("A50").Value = 'Creatina 150gr Atlas'
$ws.Range("B50").Value = 10
$ws.Range("C50").Value = 75.34
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '23943'
$ws.Range("E50").Value = 65.78
$ws.Range("F50").Value = 67.81
$ws.Range("G50").Value = 'https://www.farmaponte.com.br/creatina-150gr-atlas/p'
$ws.Range("H50").NumberFormat = "@"
$ws.Range("H50").Value = '7899732103331'
$ws.Range("I50").Value = 'Sem marca'

# Row 51
$ws.Range("A51").Value = 'Coqueteleira 1Dose Nacional Integralmedica'
$ws.Range("B51").Value = 10
$ws.Range("C51").Value = 35.56
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '24402'
$ws.Range("E51").Value = 31.04
$ws.Range("F51").Value = 32
$ws.Range("G51").Value = 'https://www.farmaponte.com.br/coqueteleira-1dose-nacional-integralmedica/p'
$ws.Range("H51").NumberFormat = "@"
$ws.Range("H51").Value = '7896311102297'
$ws.Range("I51").Value = 'Integralmédica'

# Row 52
$ws.Range("A52").Value = 'Creatina Fuel 150Gr Iridium'
$ws.Range("B52").Value = 10
$ws.Range("C52").Value = 61
$ws.Range("D52").NumberFormat = "@"
$ws.Range("D52").Value = '23969'
$ws.Range("E52").Value = 53.25
$ws.Range("F52").Value = 54.9
$ws.Range("G52").Value = 'https://www.farmaponte.com.br/creatina-fuel-150gr-iridium/p'
$ws.Range("H52").NumberFormat = "@"
$ws.Range("H52").Value = '7899732113873'
$ws.Range("I52").Value = 'Iridium'

# Row 53
$ws.Range("A53").Value = 'Creatina 90Gr Atlas'
$ws.Range("B53").Value = 10
$ws.Range("C53").Value = 56.56
$ws.Range("D53").NumberFormat = "@"
$ws.Range("D53").Value = '23966'
$ws.Range("E53").Value = 49.37
$ws.Range("F53").Value = 50.9
$ws.Range("G53").Value = 'https://www.farmaponte.com.br/creatina-90gr-atlas/p'
$ws.Range("H53").NumberFormat = "@"
$ws.Range("H53").Value = '7899732112159'
$ws.Range("I53").Value = 'Sem marca'

# Row 54
$ws.Range("A54").Value = 'Whey Wcp Concentrado Doce Leite 900Gr Iridium'
$ws.Range("B54").Value = 5
$ws.Range("C54").Value = 157.12
$ws.Range("D54").NumberFormat = "@"
$ws.Range("D54").Value = '23973'
$ws.Range("E54").Value = 145.39
$ws.Range("F54").Value = 149.89
$ws.Range("G54").Value = 'https://www.farmaponte.com.br/whey-wcp-concentrado-doce-leite-900gr-iridium/p'
$ws.Range("H54").NumberFormat = "@"
$ws.Range("H54").Value = '7899732112142'
$ws.Range("I54").Value = 'Iridium'

# Row 55
$ws.Range("A55").Value = 'Whey 100  Pote Baunilha 900GR Integralmedica'
$ws.Range("B55").Value = 10
$ws.Range("C55").Value = 184.44
$ws.Range("D55").NumberFormat = "@"
$ws.Range("D55").Value = '23836'
$ws.Range("E55").Value = 161.02
$ws.Range("F55").Value = 166
$ws.Range("G55").Value = 'https://www.farmaponte.com.br/whey-100-pote-baunilha-900gr-integralmedica/p'
$ws.Range("H55").NumberFormat = "@"
$ws.Range("H55").Value = '7896311709984'
$ws.Range("I55").Value = 'Integralmédica'

# Row 56
$ws.Range("A56").Value = 'Suplemento para Nutrição Nutren Control Diet chocolate, garrafa com 200mL'
$ws.Range("B56").Value = 0
$ws.Range("C56").Value = 14.99
$ws.Range("D56").NumberFormat = "@"
$ws.Range("D56").Value = '20792'
$ws.Range("E56").Value = 14.54
$ws.Range("F56").Value = 14.99
$ws.Range("G56").Value = 'https://www.farmaponte.com.br/suplemento-alimentar-nutren-control-chocolate-200ml/p'
$ws.Range("H56").NumberFormat = "@"
$ws.Range("H56").Value = '7891000334133'
$ws.Range("I56").Value = 'Nutren'

# Row 57
$ws.Range("A57").Value = 'Therma Pro Hardcore Com 60 Cápsulas Integralmedica'
$ws.Range("B57").Value = 10
$ws.Range("C57").Value = 73.33
$ws.Range("D57").NumberFormat = "@"
$ws.Range("D57").Value = '24798'
$ws.Range("E57").Value = 64.02
$ws.Range("F57").Value = 66
$ws.Range("G57").Value = 'https://www.farmaponte.com.br/therma-pro-hardcore-com-60-capsulas-integralmedica/p'
$ws.Range("H57").NumberFormat = "@"
$ws.Range("H57").Value = '7896311707973'
$ws.Range("I57").Value = 'Integralmédica'

# Row 58
$ws.Range("A58").Value = 'Whey 100  Pote Chocolate 900GR Integralmedica'
$ws.Range("B58").Value = 10
$ws.Range("C58").Value = 184.44
$ws.Range("D58").NumberFormat = "@"
$ws.Range("D58").Value = '23837'
$ws.Range("E58").Value = 161.02
$ws.Range("F58").Value = 166
$ws.Range("G58").Value = 'https://www.farmaponte.com.br/whey-100-pote-chocolate-900gr-integralmedica/p'
$ws.Range("H58").NumberFormat = "@"
$ws.Range("H58").Value = '7896311709991'
$ws.Range("I58").Value = 'Integralmédica'

# Row 59
$ws.Range("A59").Value = 'Fosfadil 30 cápsulas'
$ws.Range("B59").Value = 0
$ws.Range("C59").Value = 78.75
$ws.Range("D59").NumberFormat = "@"
$ws.Range("D59").Value = '18729'
$ws.Range("E59").Value = 76.39
$ws.Range("F59").Value = 78.75
$ws.Range("G59").Value = 'https://www.farmaponte.com.br/fosfadil-30cps/p'
$ws.Range("H59").NumberFormat = "@"
$ws.Range("H59").Value = '7908135001554'
$ws.Range("I59").Value = 'Sem marca'

# Row 60
$ws.Range("A60").Value = 'Kimera Thermo Iridium Labs 300mg 60 Comprimidos'
$ws.Range("B60").Value = 30
$ws.Range("C60").Value = 79.12
$ws.Range("D60").NumberFormat = "@"
$ws.Range("D60").Value = '25816'
$ws.Range("E60").NumberFormat = "@"
$ws.Range("E60").Value = ''
$ws.Range("F60").Value = 55.38
$ws.Range("G60").Value = 'https://www.farmaponte.com.br/kimera-thermo-iridium-labs-300mg-60-comprimidos/p'
$ws.Range("H60").NumberFormat = "@"
$ws.Range("H60").Value = '7899732100323'
$ws.Range("I60").Value = 'Sem marca'

# Row 61
$ws.Range("A61").Value = 'Bcaa 2400 90Caps Integralmédica'
$ws.Range("B61").Value = 10
$ws.Range("C61").Value = 76.78
$ws.Range("D61").NumberFormat = "@"
$ws.Range("D61").Value = '23793'
$ws.Range("E61").Value = 67.03
$ws.Range("F61").Value = 69.09999999999999
$ws.Range("G61").Value = 'https://www.farmaponte.com.br/bcaa-2400-90caps-integralmedica/p'
$ws.Range("H61").NumberFormat = "@"
$ws.Range("H61").Value = '7896311763269'
$ws.Range("I61").Value = 'Integralmédica'

# Row 62
$ws.Range("A62").Value = 'Nutri Whey Refil Morango 900GR Integralmedica'
$ws.Range("B62").Value = 10
$ws.Range("C62").Value = 102.43
$ws.Range("D62").NumberFormat = "@"
$ws.Range("D62").Value = '23832'
$ws.Range("E62").Value = 89.42
$ws.Range("F62").Value = 92.19
$ws.Range("G62").Value = 'https://www.farmaponte.com.br/nutri-whey-refil-morango-900gr-integralmedica/p'
$ws.Range("H62").NumberFormat = "@"
$ws.Range("H62").Value = '7896311709496'
$ws.Range("I62").Value = 'Integralmédica'

# Row 63
$ws.Range("A63").Value = 'Sinister Mass Chocolate 3KG Integralmedica'
$ws.Range("B63").Value = 10
$ws.Range("C63").Value = 148.56
$ws.Range("D63").NumberFormat = "@"
$ws.Range("D63").Value = '23834'
$ws.Range("E63").Value = 129.69
$ws.Range("F63").Value = 133.7
$ws.Range("G63").Value = 'https://www.farmaponte.com.br/mass-chocolate-3kg-integralmedica/p'
$ws.Range("H63").NumberFormat = "@"
$ws.Range("H63").Value = '7896311767298'
$ws.Range("I63").Value = 'Integralmédica'

# Row 64
$ws.Range("A64").Value = 'Whey 100  Pote Morango 900GR Integralmedica'
$ws.Range("B64").Value = 10
$ws.Range("C64").Value = 184.44
$ws.Range("D64").NumberFormat = "@"
$ws.Range("D64").Value = '23843'
$ws.Range("E64").Value = 161.02
$ws.Range("F64").Value = 166
$ws.Range("G64").Value = 'https://www.farmaponte.com.br/whey-100-pote-morango-900gr-integralmedica/p'
$ws.Range("H64").NumberFormat = "@"
$ws.Range("H64").Value = '7896311710010'
$ws.Range("I64").Value = 'Integralmédica'

# Row 65
$ws.Range("A65").Value = 'Whey Protein Shake Dux Cookies 250ml'
$ws.Range("B65").Value = 10
$ws.Range("C65").Value = 9.890000000000001
$ws.Range("D65").NumberFormat = "@"
$ws.Range("D65").Value = '24794'
$ws.Range("E65").Value = 8.630000000000001
$ws.Range("F65").Value = 8.9
$ws.Range("G65").Value = 'https://www.farmaponte.com.br/whey-protein-shake-dux-cookies-250ml/p'
$ws.Range("H65").NumberFormat = "@"
$ws.Range("H65").Value = '7898641074473'
$ws.Range("I65").Value = 'Dux Nutrition'

# Row 66
$ws.Range("A66").Value = 'Whey Protein Shake Dux Chocolate 250Ml'
$ws.Range("B66").Value = 10
$ws.Range("C66").Value = 9.890000000000001
$ws.Range("D66").NumberFormat = "@"
$ws.Range("D66").Value = '24795'
$ws.Range("E66").Value = 8.630000000000001
$ws.Range("F66").Value = 8.9
$ws.Range("G66").Value = 'https://www.farmaponte.com.br/whey-protein-shake-dux-chocolate-250ml/p'
$ws.Range("H66").NumberFormat = "@"
$ws.Range("H66").Value = '7898641074497'
$ws.Range("I66").Value = 'Dux Nutrition'

# Row 67
$ws.Range("A67").Value = 'Whey Protein Shake Dux Doce de Leite 250ml'
$ws.Range("B67").Value = 10
$ws.Range("C67").Value = 9.890000000000001
$ws.Range("D67").NumberFormat = "@"
$ws.Range("D67").Value = '24792'
$ws.Range("E67").Value = 8.630000000000001
$ws.Range("F67").Value = 8.9
$ws.Range("G67").Value = 'https://www.farmaponte.com.br/whey-protein-shake-dux-doce-de-leite-250ml/p'
$ws.Range("H67").NumberFormat = "@"
$ws.Range("H67").Value = '7898641074480'
$ws.Range("I67").Value = 'Dux Nutrition'

# Row 68
$ws.Range("A68").Value = 'Suplemento Alimentar Florence Pro 6 Sachês de 4g cada'
$ws.Range("B68").Value = 30
$ws.Range("C68").Value = 42.19
$ws.Range("D68").NumberFormat = "@"
$ws.Range("D68").Value = '23798'
$ws.Range("E68").Value = 28.52
$ws.Range("F68").Value = 29.4
$ws.Range("G68").Value = 'https://www.farmaponte.com.br/suplemento-alimentar-florence-pro-6-saches-de-4g-cada-biolab/p'
$ws.Range("H68").NumberFormat = "@"
$ws.Range("H68").Value = '7896112407201'
$ws.Range("I68").Value = 'Sem marca'

# Row 69
$ws.Range("A69").Value = 'Sinister Mass Baunilha 3KG Integralmedica'
$ws.Range("B69").Value = 10
$ws.Range("C69").Value = 148.56
$ws.Range("D69").NumberFormat = "@"
$ws.Range("D69").Value = '23833'
$ws.Range("E69").Value = 129.69
$ws.Range("F69").Value = 133.7
$ws.Range("G69").Value = 'https://www.farmaponte.com.br/mass-baunilha-3kg-integralmedica/p'
$ws.Range("H69").NumberFormat = "@"
$ws.Range("H69").Value = '7896311767281'
$ws.Range("I69").Value = 'Integralmédica'

# Row 70
$ws.Range("A70").Value = 'Whey Protein Dux Chocolate Branco 250ml'
$ws.Range("B70").Value = 10
$ws.Range("C70").Value = 9.890000000000001
$ws.Range("D70").NumberFormat = "@"
$ws.Range("D70").Value = '24793'
$ws.Range("E70").Value = 8.630000000000001
$ws.Range("F70").Value = 8.9
$ws.Range("G70").Value = 'https://www.farmaponte.com.br/whey-protein-dux-chocolate-branco-250ml/p'
$ws.Range("H70").NumberFormat = "@"
$ws.Range("H70").Value = '7898641074503'
$ws.Range("I70").Value = 'Dux Nutrition'

# Row 71
$ws.Range("A71").Value = 'Isotônico Powerade Limão 500ml'
$ws.Range("B71").Value = 10
$ws.Range("C71").Value = 6.99
$ws.Range("D71").NumberFormat = "@"
$ws.Range("D71").Value = '22732'
$ws.Range("E71").Value = 6.1
$ws.Range("F71").Value = 6.29
$ws.Range("G71").Value = 'https://www.farmaponte.com.br/powerade-limao-500ml/p'
$ws.Range("H71").NumberFormat = "@"
$ws.Range("H71").Value = '7894900500035'
$ws.Range("I71").Value = 'Powerade'

# Row 72
$ws.Range("A72").Value = 'Isotônico Powerade Mountain Blast 500ml'
$ws.Range("B72").Value = 10
$ws.Range("C72").Value = 6.99
$ws.Range("D72").NumberFormat = "@"
$ws.Range("D72").Value = '24783'
$ws.Range("E72").Value = 6.1
$ws.Range("F72").Value = 6.29
$ws.Range("G72").Value = 'https://www.farmaponte.com.br/isotonico-powerade-mountain-blast-500ml/p'
$ws.Range("H72").NumberFormat = "@"
$ws.Range("H72").Value = '7894900504002'
$ws.Range("I72").Value = 'Sem marca'

# Row 73
$ws.Range("A73").Value = 'Isotônico Powerade Sabor Frutas Tropicais 500ml'
$ws.Range("B73").Value = 10
$ws.Range("C73").Value = 6.99
$ws.Range("D73").NumberFormat = "@"
$ws.Range("D73").Value = '24813'
$ws.Range("E73").Value = 6.1
$ws.Range("F73").Value = 6.29
$ws.Range("G73").Value = 'https://www.farmaponte.com.br/isotonico-powerade-sabor-frutas-tropicais-500ml/p'
$ws.Range("H73").NumberFormat = "@"
$ws.Range("H73").Value = '7894900508017'
$ws.Range("I73").Value = 'Sem marca'

# Row 74
$ws.Range("A74").Value = 'Go Energy Now Gel Guarana Com Acai 30gr'
$ws.Range("B74").Value = 10
$ws.Range("C74").Value = 4.82
$ws.Range("D74").NumberFormat = "@"
$ws.Range("D74").Value = '25480'
$ws.Range("E74").Value = 4.21
$ws.Range("F74").Value = 4.34
$ws.Range("G74").Value = 'https://www.farmaponte.com.br/go-energy-now-gel-guarana-com-acai-30gr/p'
$ws.Range("H74").NumberFormat = "@"
$ws.Range("H74").Value = '7899621106573'
$ws.Range("I74").Value = 'Sem marca'
